# Report template formatting tweaks
# Adjust the paragraph indentation for the Heading1/Heading2/Heading3
# styles that drive the multilevel "report" numbering (numId 17 ->
# abstractNumId 3).
#
# Word stores indentation in twips (1/20 pt) in the OOXML, but the
# ParagraphFormat object works in points, so divide the target twip
# values by 20 before assigning them.

$d = $word.ActiveDocument

# Heading 1: add a first-line indent of 267 twips (13.35 pt).
$h1 = $d.Styles("Heading1")
$h1.ParagraphFormat.FirstLineIndent = 267 / 20

# Heading 2: left indent of 567 twips with a matching hanging indent
# (i.e. FirstLineIndent = -567 twips) so the wrapped lines line up
# under the numbering.
$h2 = $d.Styles("Heading2")
$h2.ParagraphFormat.LeftIndent = 567 / 20
$h2.ParagraphFormat.FirstLineIndent = -567 / 20

# Heading 3: left indent of 1276 twips with a 142 twip hanging indent.
$h3 = $d.Styles("Heading3")
$h3.ParagraphFormat.LeftIndent = 1276 / 20
$h3.ParagraphFormat.FirstLineIndent = -142 / 20
